# edit.ps1 - apply the commit's changes to the loaded document
# 1. Address line: "MALKLOAD ROAD,LHR ZONE" -> "CHUMA,LHR ZONE"
# 2. Table row (Monthly general Servicing): description, Qty, Rate, Amount

$d = $word.ActiveDocument

# Helper: replace text strictly within a given Range's [Start, End) bounds.
# Using Find.Execute with Replace = wdReplaceAll (2) on a sub-range was
# observed to search/replace across the *entire* document story rather than
# staying confined to the range, so we use wdReplaceOne (1) on a range
# rebuilt from explicit Start/End offsets, which stays correctly scoped.
function Replace-InRange($range, [string]$oldText, [string]$newText) {
    $scoped = $d.Range($range.Start, $range.End)
    $scoped.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                          $true, 0, $false, $newText, 1) | Out-Null
}

# --- 1. Fix the address line ---
Replace-InRange $d.Content "MALKLOAD ROAD,LHR ZONE" "CHUMA,LHR ZONE"

# --- 2. Update the service row in the pricing table ---
$table = $d.Tables.Item(1)

# Row 3: "Monthly general Servicing (upto to 2 Ton) " | Qty 2 | Rate 300.0 | Amount 600.0
Replace-InRange $table.Cell(3, 2).Range "Monthly general Servicing (upto to 2 Ton) " "Gas flushing "
Replace-InRange $table.Cell(3, 3).Range "2" "3"
Replace-InRange $table.Cell(3, 4).Range "300.0" "0.0"
Replace-InRange $table.Cell(3, 5).Range "600.0" "0.0"
